$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new tag columns right after the "keyvault" column (V), i.e. before the
# old "payment-details" column (W). This shifts every column from W onward five
# places to the right (W->AB, X->AC, ... AS->AX), matching the diff exactly.
$ws.Range("W1:AA1").EntireColumn.Insert()

# Populate the header row for the 5 newly inserted columns.
$ws.Range("W1").Value = "microsoft.app/containerapps"
$ws.Range("X1").Value = "microsoft.cache/redis"
$ws.Range("Y1").Value = "microsoft.keyvault/vaults"
$ws.Range("Z1").Value = "microsoft.storage/storageaccounts"
$ws.Range("AA1").Value = "microsoft.web/sites"

# Mark the relevant elements with the newly introduced Azure-resource-type tags.
$ws.Range("W2").Value = "X"   # s184d01-comp-complete-app        -> microsoft.app/containerapps
$ws.Range("W3").Value = "X"   # s184d01-comp-complete-app-worker -> microsoft.app/containerapps
$ws.Range("Y4").Value = "X"   # s184d01-comp-tfvars              -> microsoft.keyvault/vaults
$ws.Range("X5").Value = "X"   # s184d01-compdefault              -> microsoft.cache/redis
$ws.Range("Y6").Value = "X"   # ssphp-metrics                    -> microsoft.keyvault/vaults
$ws.Range("AA7").Value = "X"  # ssphp-metrics-rust-p3sha         -> microsoft.web/sites
$ws.Range("Z8").Value = "X"   # tfstatel95cd                     -> microsoft.storage/storageaccounts
$ws.Range("Z9").Value = "X"   # tfstatep3sha                     -> microsoft.storage/storageaccounts
